# Removed Test Case Inter-Dependency
# Rename the product name string (shared by both sheets), give the
# product a distinct "shortname" text token ("413w") instead of the
# numeric product id, and switch the active sheet/tab selection from
# ProductLoanInput to ProductLoanOutput.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Rename the product name string used across the workbook. Both cells
# reference the same shared string, so updating both to the same new
# text edits that shared string entry in place.
$wsInput.Range("B1").Value = "4131-RBI-EI-FL-DL-NOREC-MOREREPAY-1st"
$wsOutput.Range("B1").Value = "4131-RBI-EI-FL-DL-NOREC-MOREREPAY-1st"

# Update B2 (shortname) to a new text value, decoupling it from the
# numeric product id used elsewhere (removes inter-test dependency).
$wsInput.Range("B2").Value = "413w"

# Make ProductLoanOutput the active/selected sheet (tabSelected moves
# from sheet1 to sheet2).
$wsOutput.Activate()
